# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 1 de Mayo de 2020 a las 15:22"

# Estados Unidos (row 4)
$ws.Range("B4").Value = 1095977
$ws.Range("C4").Value = 954
$ws.Range("E4").Value = 876364
$ws.Range("G4").Value = 20
$ws.Range("H4").Value = 63876

# Arabia Saudita (row 23)
$ws.Range("B23").Value = 24097
$ws.Range("C23").Value = 1344
$ws.Range("D23").Value = 3555
$ws.Range("E23").Value = 20373
$ws.Range("G23").Value = 7
$ws.Range("H23").Value = 169

# Serbia (row 42)
$ws.Range("B42").Value = 9205
$ws.Range("C42").Value = 196
$ws.Range("D42").Value = 1379
$ws.Range("E42").Value = 7641
$ws.Range("G42").Value = 6
$ws.Range("H42").Value = 185

# Islandia (row 76)
$ws.Range("B76").Value = 1798
$ws.Range("C76").Value = 1
$ws.Range("D76").Value = 1689
$ws.Range("E76").Value = 99

# Sri Lanka (row 103)
$ws.Range("B103").Value = 674
$ws.Range("C103").Value = 11
$ws.Range("E103").Value = 510
